$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "SamplesTab" row (row 3) keeps its query in column B ("TabQuery"). The
# regression-suite query is trimmed down: it no longer pulls the tumor
# status / analyte type columns out of df_sample.
$newQuery = @"
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession 
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs002529'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
"@

$ws.Range("B3").Value = $newQuery

# The sheet had scrolled so column B was pinned at the left edge
# (topLeftCell = B3); scroll it back so column A is visible again while
# keeping the same active selection (C3).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("C3").Select()
